$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    # Force the literal text into the cell without Excel re-interpreting
    # numeric-looking strings (e.g. "1.00", "9.50") as numbers, which would
    # silently drop significant trailing zeros / introduce float noise.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextCell "D2" "65.914.91"
$ws.Range("E2").Value = "  +0.98%  "
Set-TextCell "D3" "3.309.85"
$ws.Range("E3").Value = "  +0.36%  "
Set-TextCell "D4" "1.00"
$ws.Range("E4").Value = "  -0.14%  "
Set-TextCell "D5" "187.46"
$ws.Range("E5").Value = "  +4.76%  "
Set-TextCell "D6" "556.14"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -1.07%  "
Set-TextCell "D9" "3.301.03"
$ws.Range("E9").Value = "  +0.19%  "
Set-TextCell "D10" "0.181"
$ws.Range("E10").Value = "  -1.90%  "
Set-TextCell "D11" "0.581"
$ws.Range("E11").Value = "  -0.26%  "
Set-TextCell "D12" "46.83"
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("E14").Value = "  +1.32%  "
Set-TextCell "D15" "3.838.06"
$ws.Range("E15").Value = "  -0.30%  "
Set-TextCell "D16" "601.45"
$ws.Range("E16").Value = "  +0.52%  "
Set-TextCell "D17" "65.882.28"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell "D18" "17.89"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell "D19" "0.117"
$ws.Range("E19").Value = "  +0.82%  "
Set-TextCell "D20" "3.293.47"
$ws.Range("E20").Value = "  -0.31%  "
Set-TextCell "D21" "11.04"
$ws.Range("E21").Value = "  -2.62%  "
Set-TextCell "D22" "0.899"
$ws.Range("E22").Value = "  +0.16%  "
Set-TextCell "D23" "18.64"
$ws.Range("E23").Value = "  +7.63%  "
Set-TextCell "D24" "5.06"
$ws.Range("E24").Value = "  +1.07%  "
Set-TextCell "D25" "100.74"
$ws.Range("E25").Value = "  -1.48%  "
Set-TextCell "D26" "3.94"
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("E27").Value = "  +3.21%  "
$ws.Range("E28").Value = "  -1.16%  "
Set-TextCell "D29" "9.50"
$ws.Range("E29").Value = "  +2.18%  "
Set-TextCell "D30" "8.67"
$ws.Range("E30").Value = "  +0.80%  "
Set-TextCell "D31" "30.23"
$ws.Range("E31").Value = "  -0.26%  "
Set-TextCell "D32" "6.72"
$ws.Range("E32").Value = "  +7.85%  "
Set-TextCell "D33" "3.86"
$ws.Range("E33").Value = "  +0.31%  "
Set-TextCell "D34" "572.14"
$ws.Range("E34").Value = "  +8.07%  "
Set-TextCell "D35" "11.01"
$ws.Range("E35").Value = "  +0.39%  "
$ws.Range("E36").Value = "  +0.31%  "
Set-TextCell "D37" "1.00"
$ws.Range("E37").Value = "  +0.07%  "
Set-TextCell "D38" "3.697.61"
$ws.Range("E38").Value = "  -2.91%  "
Set-TextCell "D39" "56.84"
$ws.Range("E39").Value = "  +1.42%  "
Set-TextCell "D40" "33.73"
$ws.Range("E40").Value = "  +6.65%  "
$ws.Range("E41").Value = "  +9.15%  "
Set-TextCell "D42" "0.129"
$ws.Range("E42").Value = "  +2.83%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D43" "3.25"
$ws.Range("E43").Value = "  -5.12%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell "D44" "0.0₃0706"
$ws.Range("E44").Value = "  +0.47%  "
Set-TextCell "D45" "2.65"
$ws.Range("E45").Value = "  +0.45%  "
Set-TextCell "D46" "3.38"
$ws.Range("E46").Value = "  +4.42%  "
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("E48").Value = "  +2.45%  "
$ws.Range("E49").Value = "  -0.11%  "
Set-TextCell "D50" "2.56"
$ws.Range("E50").Value = "  -0.56%  "
Set-TextCell "D51" "0.999"
$ws.Range("E51").Value = "  -0.21%  "
